# Fruta / hortaliza, semanal
# Insert one new weekly record as row 181 (pushing the existing rows 181-296
# down to 182-297) on the single worksheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 181..296 down one row, creating a fresh (blank) row 181.
$ws.Rows.Item(181).Insert()

# Populate the newly inserted row 181 with the new weekly reading.
$ws.Range("A181").Value2 = 1
$ws.Range("B181").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C181").Value2 = "Arica y Parinacota"
$ws.Range("D181").Value2 = 44824
$ws.Range("E181").Value2 = 15
$ws.Range("F181").Value2 = "Fruta"
$ws.Range("G181").Value2 = 100108
$ws.Range("H181").Value2 = "Tropicales y subtropicales"
$ws.Range("I181").Value2 = 100108006
$ws.Range("J181").Value2 = "Plátano"
$ws.Range("K181").Value2 = "Sin especificar"
$ws.Range("L181").Value2 = "Pintón"
$ws.Range("M181").Value2 = 120
$ws.Range("N181").Value2 = 26000
$ws.Range("O181").Value2 = 27000
$ws.Range("P181").Value2 = 26500
$ws.Range("Q181").Value2 = "$/caja 20 kilos"
$ws.Range("R181").Value2 = "Ecuador"
$ws.Range("S181").Value2 = 1325
$ws.Range("T181").Value2 = 20
